$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44, "test test", "04-11-2023"),
    @(45, "test again", "04-11-2023"),
    @(46, "testestes", "04-11-2023"),
    @(47, "tstest", "04-11-2023")
)

$row = 46
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
